# Append 20 new "Math-Solver" test entries (rows 43-62) to Tabelle1 (sheet1),
# matching the new image/formula/result samples added upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
  @(40, "(13+11)/2", 12),
  @(41, "(14*3)-1", 41),
  @(42, "((7+4)*30)/4", 82.5),
  @(43, "500-40", 460),
  @(44, "1337-27", 1300),
  @(45, "(42+34)*9", 684),
  @(46, "(47*11)/(4-2)", 258.5),
  @(47, "(0-9)/(-7)", 1.2857),
  @(48, "((679+1)/2)*4", 1360),
  @(49, "(49+7-4)*5", 260),
  @(50, "717+471/-7", 649.714),
  @(51, "874+37*111", 4981),
  @(52, "(69+69)*71", 9798),
  @(53, "666-111*3", 333),
  @(54, "(1995-21)*3", 5922),
  @(55, "(50/5)*(5-4)", 10),
  @(56, "(612-4+7)/2", 307.5),
  @(57, "46-0+8", 54),
  @(58, "8+8", 16),
  @(59, "7-7", 0)
)

$startRow = 43
$r = $startRow
foreach ($item in $data) {
  $ws.Cells.Item($r, 1).Value = $item[0]
  $ws.Cells.Item($r, 2).Value = $item[1]
  $ws.Cells.Item($r, 3).Value = $item[2]
  $r++
}
$endRow = $r - 1

# Visual formatting for the new block: thin "card" look with a light-gray
# medium border, wrapped text, numbers/results right-aligned.
$colA = $ws.Range("A" + $startRow + ":A" + $endRow)
$colA.WrapText = $true
$colA.HorizontalAlignment = -4152
$colA.Borders.Weight = -4138
$colA.Borders.Color = 13421772

$colC = $ws.Range("C" + $startRow + ":C" + $endRow)
$colC.WrapText = $true
$colC.HorizontalAlignment = -4152
$colC.Borders.Weight = -4138
$colC.Borders.Color = 13421772

$formulaCol = $ws.Range("B" + $startRow + ":B" + $endRow)
$formulaCol.WrapText = $true
$formulaCol.Borders.Weight = -4138
$formulaCol.Borders.Color = 13421772

# Last new row's formula cell ("7-7") was kept as explicit text.
$lastFormulaCell = $ws.Range("B" + $endRow)
$lastFormulaCell.NumberFormat = "@"
$lastFormulaCell.HorizontalAlignment = -4152
